# Add new columns I (I0) and J (IF) to Sheet1, mirroring the existing
# header/data layout used by columns A-H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from an existing
# header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-31) ---------------------------------------------
$iValues = @(8,5,6,5,8,5,5,6,6,9,8,8,5,9,9,5,9,7,5,7,9,6,10,7,6,8,5,8,9,5)
$jValues = @(8,6,6,5,9,5,6,7,7,9,8,8,7,9,9,5,9,7,6,9,9,7,11,8,6,8,6,8,9,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

# --- Update the sheet dimension to reflect the new columns --------------
$ws.Range("A1:J31").Select() | Out-Null

$wb.Save()
